# Seguimiento et3.xlsx - revisiones del recurso 5 (fin de primera semana)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Tarea 19 (row 22): estado pasa de "Sin empezar" a "Finalizada/Corregida" ---
$ws.Range("C22").Value = "Finalizada/Corregida"
$ws.Rows.Item(22).RowHeight = 21.65

# --- Tarea 22 (row 25): estado + observaciones + recurso asociado revisados ---
$ws.Range("C25").Value = "Finalizada/Corregida"
$ws.Range("D25").Value = "Corregidas algunas referencias a controladores y clases. Eliminado edit porque los atributos son clave."
$ws.Range("G25").Value = 20
$ws.Rows.Item(25).RowHeight = 52.4

# --- Tarea 23 (row 26): estado pasa de "Pendiente de correción" a "Finalizada/Corregida" ---
$ws.Range("C26").Value = "Finalizada/Corregida"
$ws.Rows.Item(26).RowHeight = 31.9

# --- Vista de la hoja: nueva celda activa y desplazamiento al inicio ---
$aw = $excel.ActiveWindow
$aw.ScrollRow = 1
$aw.ScrollColumn = 1
$ws.Range("G16").Select() | Out-Null
